$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17 (hunk 0)
$ws.Cells.Item(17, 8).Value = 1112.5
$ws.Cells.Item(17, 10).Value = 1112.5
$ws.Cells.Item(17, 12).Value = 3337.5
$ws.Cells.Item(17, 14).Value = -3673.5

# Row 33 (hunk 1)
$ws.Cells.Item(33, 8).Value = 331.57693
$ws.Cells.Item(33, 9).Value = 328.61905
$ws.Cells.Item(33, 10).Value = 344
$ws.Cells.Item(33, 11).Value = 328.61905
$ws.Cells.Item(33, 12).Value = 344
$ws.Cells.Item(33, 13).Value = -99.61905000000002
$ws.Cells.Item(33, 14).Value = -802

# Row 58 (hunk 2)
$ws.Cells.Item(58, 8).Value = 733.5238000000001
$ws.Cells.Item(58, 9).Value = 380.5
$ws.Cells.Item(58, 10).Value = 1439.5714
$ws.Cells.Item(58, 11).Value = 1141.5
$ws.Cells.Item(58, 12).Value = 4318.7142
$ws.Cells.Item(58, 13).Value = -991.5
$ws.Cells.Item(58, 14).Value = -4618.7142

# Row 74 (hunk 3)
$ws.Cells.Item(74, 8).Value = 3633
$ws.Cells.Item(74, 9).Value = 3700
$ws.Cells.Item(74, 10).Value = 3499
$ws.Cells.Item(74, 11).Value = 3700
$ws.Cells.Item(74, 12).Value = 3499
$ws.Cells.Item(74, 13).Value = -2764
$ws.Cells.Item(74, 14).Value = -5371

# Row 77 (hunk 4)
$ws.Cells.Item(77, 8).Value = 3633
$ws.Cells.Item(77, 9).Value = 3700
$ws.Cells.Item(77, 10).Value = 3499
$ws.Cells.Item(77, 11).Value = 18500
$ws.Cells.Item(77, 12).Value = 17495
$ws.Cells.Item(77, 13).Value = -13820
$ws.Cells.Item(77, 14).Value = -26855

# Row 112 (hunk 5)
$ws.Cells.Item(112, 8).Value = 2916.2354
$ws.Cells.Item(112, 10).Value = 3029.75
$ws.Cells.Item(112, 12).Value = 9089.25
$ws.Cells.Item(112, 14).Value = -11305.25

# Row 132 (hunk 6)
$ws.Cells.Item(132, 8).Value = 9528831
$ws.Cells.Item(132, 9).Value = 10758092
$ws.Cells.Item(132, 11).Value = 32274276
$ws.Cells.Item(132, 13).Value = -32271746

# Row 138 (hunk 7)
$ws.Cells.Item(138, 8).Value = 1478.1398
$ws.Cells.Item(138, 9).Value = 868.7917
$ws.Cells.Item(138, 10).Value = 1690.0869
$ws.Cells.Item(138, 11).Value = 2606.3751
$ws.Cells.Item(138, 12).Value = 5070.2607
$ws.Cells.Item(138, 13).Value = 2533.6249
$ws.Cells.Item(138, 14).Value = -15350.2607

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (hunk 8)
$ws.Cells.Item(32, 8).Value = 8898.940000000001
$ws.Cells.Item(32, 9).Value = 6704.488
$ws.Cells.Item(32, 11).Value = 6704.488
$ws.Cells.Item(32, 13).Value = -6417.488

# Row 139 (hunk 9)
$ws.Cells.Item(139, 8).Value = 31445
$ws.Cells.Item(139, 10).Value = 31445
$ws.Cells.Item(139, 12).Value = 31445
$ws.Cells.Item(139, 14).Value = -41725

$ws = $wb.Worksheets.Item("BSM")
# Row 94 (hunk 10)
$ws.Cells.Item(94, 8).Value = 20834200
$ws.Cells.Item(94, 9).Value = 31250302
$ws.Cells.Item(94, 10).Value = 1996.25
$ws.Cells.Item(94, 11).Value = 31250302
$ws.Cells.Item(94, 12).Value = 1996.25
$ws.Cells.Item(94, 13).Value = -31249851
$ws.Cells.Item(94, 14).Value = -2898.25

# Row 107 (hunk 11)
$ws.Cells.Item(107, 8).Value = 2462.1428
$ws.Cells.Item(107, 9).Value = 1984.4
$ws.Cells.Item(107, 11).Value = 1984.4
$ws.Cells.Item(107, 13).Value = -64.40000000000009

# Row 134 (hunk 12)
$ws.Cells.Item(134, 8).Value = 1855.5186
$ws.Cells.Item(134, 9).Value = 1406.875
$ws.Cells.Item(134, 10).Value = 5444.6665
$ws.Cells.Item(134, 11).Value = 4220.625
$ws.Cells.Item(134, 12).Value = 16333.9995
$ws.Cells.Item(134, 13).Value = -1685.625
$ws.Cells.Item(134, 14).Value = -21403.9995

$ws = $wb.Worksheets.Item("CRP")
# Row 16 (hunk 13)
$ws.Cells.Item(16, 8).Value = 47620416
$ws.Cells.Item(16, 9).Value = 58824920
$ws.Cells.Item(16, 10).Value = 1275
$ws.Cells.Item(16, 11).Value = 58824920
$ws.Cells.Item(16, 12).Value = 1275
$ws.Cells.Item(16, 13).Value = -58824633
$ws.Cells.Item(16, 14).Value = -1849

# Row 31 (hunk 14)
$ws.Cells.Item(31, 8).Value = 1363.25
$ws.Cells.Item(31, 9).Value = 1119.625
$ws.Cells.Item(31, 10).Value = 2825
$ws.Cells.Item(31, 11).Value = 1119.625
$ws.Cells.Item(31, 12).Value = 2825
$ws.Cells.Item(31, 13).Value = -824.625
$ws.Cells.Item(31, 14).Value = -3415

# Row 34 (hunk 15)
$ws.Cells.Item(34, 8).Value = 1363.25
$ws.Cells.Item(34, 9).Value = 1119.625
$ws.Cells.Item(34, 10).Value = 2825
$ws.Cells.Item(34, 11).Value = 1119.625
$ws.Cells.Item(34, 12).Value = 2825
$ws.Cells.Item(34, 13).Value = -917.625
$ws.Cells.Item(34, 14).Value = -3229

# Row 58 (hunk 16)
$ws.Cells.Item(58, 8).Value = 1567.5834
$ws.Cells.Item(58, 10).Value = 2202.4
$ws.Cells.Item(58, 12).Value = 2202.4
$ws.Cells.Item(58, 14).Value = -2608.4

# Row 99 (hunk 17)
$ws.Cells.Item(99, 8).Value = 1590.8462
$ws.Cells.Item(99, 9).Value = 1468.1
$ws.Cells.Item(99, 11).Value = 1468.1
$ws.Cells.Item(99, 13).Value = 29.90000000000009

# Row 113 (hunk 18)
$ws.Cells.Item(113, 8).Value = 47620416
$ws.Cells.Item(113, 9).Value = 58824920
$ws.Cells.Item(113, 10).Value = 1275
$ws.Cells.Item(113, 11).Value = 58824920
$ws.Cells.Item(113, 12).Value = 1275
$ws.Cells.Item(113, 13).Value = -58822750
$ws.Cells.Item(113, 14).Value = -5615

# Row 126 (hunk 19)
$ws.Cells.Item(126, 8).Value = 1590.8462
$ws.Cells.Item(126, 9).Value = 1468.1
$ws.Cells.Item(126, 11).Value = 4404.299999999999
$ws.Cells.Item(126, 13).Value = -1934.299999999999

# Row 136 (hunk 20)
$ws.Cells.Item(136, 8).Value = 1567.5834
$ws.Cells.Item(136, 10).Value = 2202.4
$ws.Cells.Item(136, 12).Value = 6607.200000000001
$ws.Cells.Item(136, 14).Value = -11707.2

$ws = $wb.Worksheets.Item("CUL")
# Row 69 (hunk 21)
$ws.Cells.Item(69, 8).Value = 3652.3333
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 3652.3333
$ws.Cells.Item(69, 11).Value = 0
$ws.Cells.Item(69, 12).ClearContents()
$ws.Cells.Item(69, 13).Value = 10956.9999
$ws.Cells.Item(69, 14).Value = -12578.9999

# Row 72 (hunk 22)
$ws.Cells.Item(72, 8).Value = 3652.3333
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 3652.3333
$ws.Cells.Item(72, 11).Value = 0
$ws.Cells.Item(72, 12).ClearContents()
$ws.Cells.Item(72, 13).Value = 32870.9997
$ws.Cells.Item(72, 14).Value = -40982.9997

# Row 98 (hunk 23)
$ws.Cells.Item(98, 8).Value = 1668.125
$ws.Cells.Item(98, 9).Value = 2066.8333
$ws.Cells.Item(98, 10).Value = 472
$ws.Cells.Item(98, 11).Value = 6200.499899999999
$ws.Cells.Item(98, 12).Value = 1416
$ws.Cells.Item(98, 13).Value = -4702.499899999999
$ws.Cells.Item(98, 14).Value = -4412

# Row 122 (hunk 24)
$ws.Cells.Item(122, 8).Value = 1067.9048
$ws.Cells.Item(122, 10).Value = 1163.8667
$ws.Cells.Item(122, 12).Value = 10474.8003
$ws.Cells.Item(122, 14).Value = -15374.8003

# Row 131 (hunk 25)
$ws.Cells.Item(131, 8).Value = 23259628
$ws.Cells.Item(131, 10).Value = 4386.6484
$ws.Cells.Item(131, 12).Value = 13159.9452
$ws.Cells.Item(131, 14).Value = -23239.9452

# Row 134 (hunk 26)
$ws.Cells.Item(134, 8).Value = 3314.3333
$ws.Cells.Item(134, 9).Value = 1624.5883
$ws.Cells.Item(134, 10).Value = 5109.6875
$ws.Cells.Item(134, 11).Value = 4873.7649
$ws.Cells.Item(134, 12).Value = 15329.0625
$ws.Cells.Item(134, 13).Value = 196.2350999999999
$ws.Cells.Item(134, 14).Value = -25469.0625

# Row 136 (hunk 27)
$ws.Cells.Item(136, 8).Value = 1875.3889
$ws.Cells.Item(136, 9).Value = 1171.5834
$ws.Cells.Item(136, 10).Value = 3283
$ws.Cells.Item(136, 11).Value = 3514.7502
$ws.Cells.Item(136, 12).Value = 9849
$ws.Cells.Item(136, 13).Value = 1585.2498
$ws.Cells.Item(136, 14).Value = -20049

# Row 138 (hunk 28)
$ws.Cells.Item(138, 8).Value = 1894.4231
$ws.Cells.Item(138, 10).Value = 2242.7
$ws.Cells.Item(138, 12).Value = 6728.099999999999
$ws.Cells.Item(138, 14).Value = -17008.1

# Row 139 (hunk 29)
$ws.Cells.Item(139, 8).Value = 1661.2565
$ws.Cells.Item(139, 9).Value = 1733.2778
$ws.Cells.Item(139, 10).Value = 1599.5238
$ws.Cells.Item(139, 11).Value = 5199.8334
$ws.Cells.Item(139, 12).Value = 4798.5714
$ws.Cells.Item(139, 13).Value = -59.83340000000044
$ws.Cells.Item(139, 14).Value = -15078.5714

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (hunk 30)
$ws.Cells.Item(2, 8).Value = 173.35294
$ws.Cells.Item(2, 9).Value = 155.625
$ws.Cells.Item(2, 10).Value = 189.11111
$ws.Cells.Item(2, 11).Value = 155.625
$ws.Cells.Item(2, 12).Value = 189.11111
$ws.Cells.Item(2, 13).Value = -42.625
$ws.Cells.Item(2, 14).Value = -415.11111

# Row 102 (hunk 31)
$ws.Cells.Item(102, 8).Value = 1222.0667
$ws.Cells.Item(102, 9).Value = 994.6923
$ws.Cells.Item(102, 10).Value = 2700
$ws.Cells.Item(102, 11).Value = 994.6923
$ws.Cells.Item(102, 12).Value = 2700
$ws.Cells.Item(102, 13).Value = 627.3077
$ws.Cells.Item(102, 14).Value = -5944

# Row 113 (hunk 32)
$ws.Cells.Item(113, 8).Value = 2014.2
$ws.Cells.Item(113, 10).Value = 2500
$ws.Cells.Item(113, 12).Value = 2500
$ws.Cells.Item(113, 14).Value = -6840

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (hunk 33)
$ws.Cells.Item(40, 8).Value = 2433.8696
$ws.Cells.Item(40, 9).Value = 1740.5
$ws.Cells.Item(40, 10).Value = 4930
$ws.Cells.Item(40, 11).Value = 1740.5
$ws.Cells.Item(40, 12).Value = 4930
$ws.Cells.Item(40, 13).Value = -1604.5
$ws.Cells.Item(40, 14).Value = -5202

# Row 93 (hunk 34)
$ws.Cells.Item(93, 8).Value = 954
$ws.Cells.Item(93, 9).Value = 903
$ws.Cells.Item(93, 10).Value = 979.5
$ws.Cells.Item(93, 11).Value = 903
$ws.Cells.Item(93, 12).Value = 979.5
$ws.Cells.Item(93, 13).Value = 345
$ws.Cells.Item(93, 14).Value = -3475.5

# Row 134 (hunk 35)
$ws.Cells.Item(134, 8).Value = 32294.285
$ws.Cells.Item(134, 10).Value = 32294.285
$ws.Cells.Item(134, 12).Value = 32294.285
$ws.Cells.Item(134, 14).Value = -42434.285

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (hunk 36)
$ws.Cells.Item(132, 8).Value = 3770.5789
$ws.Cells.Item(132, 9).Value = 4091.8572
$ws.Cells.Item(132, 11).Value = 12275.5716
$ws.Cells.Item(132, 13).Value = -9745.571599999999

# Row 140 (hunk 37)
$ws.Cells.Item(140, 8).Value = 32609.8
$ws.Cells.Item(140, 10).Value = 32609.8
$ws.Cells.Item(140, 12).Value = 32609.8
$ws.Cells.Item(140, 14).Value = -42969.8
